$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "204264543"
$ws.Range("B3").Value = "Dan"
$ws.Range("C3").Value = "marinescu"
$ws.Range("D3").Value = "1234"
$ws.Range("E3").Value = $true

$ws.Range("A3").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
